$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Trad Room", 4.3, 168, "Japanese,Sushi Bars,Cocktail Bars"),
    @("Tradisyon", 4.3, 284, "Filipino"),
    @("Dolly Varden", 4.5, 531, "Cocktail Bars,New American,Breakfast & Brunch"),
    @("Friend Of A Farmer - Upper West Side", 4.3, 35, "American,Comfort Food,Breakfast & Brunch,Upper West Side"),
    @("Winona’s", 4.5, 53, "American,Bedford Stuyvesant"),
    @("Butler", 4.2, 209, "Coffee & Tea,Cafes"),
    @("Carla", 4.4, 176, "New American,Burgers,Cocktail Bars"),
    @("Bolivian Llama Party", 4.7, 212, "Latin American,Sunnyside"),
    @("5ive Spice", 4.6, 749, "Vietnamese,Tacos,Sandwiches"),
    @("Grindhaus", 4.5, 116, "New American")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
